$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "bandwidth" columns (cm-bandwidth, icm8-bandwidth, icm16-bandwidth,
# icm32-bandwidth) occupied columns C:F. The speedup columns (cm-speedup,
# icm8-speedup, icm16-speedup, icm32-speedup) occupied columns G:J and held
# the important information that should be kept. Per the commit message,
# we delete the bandwidth columns entirely; the speedup columns then shift
# left to become the new C:F.
$deleteRange = $ws.Range("C1:F1048576")
$deleteRange.Select() | Out-Null
$deleteRange.EntireColumn.Delete() | Out-Null

# Reflect the resulting selection state (columns C:F now hold the former
# speedup data), matching the saved selection in the workbook.
$ws.Range("C1:F1048576").Select() | Out-Null
$excel.ActiveCell = $ws.Range("C1")
